# Daily attendance processing - 2025-11-07 08:53:51
# Normalizes the "Recorded By" column (G) so that entries where "System"
# is not already the first listed recorder have their two comma-separated
# values swapped (e.g. "user@example.com, System" -> "System, user@example.com").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the last used row based on the worksheet's used range.
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

# Column G holds "Recorded By". Data starts on row 2 (row 1 is the header).
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($null -eq $val) { continue }

    $text = [string]$val
    $parts = $text -split ", "

    if ($parts.Length -eq 2 -and $parts[0] -ne "System") {
        $newText = $parts[1] + ", " + $parts[0]
        $cell.Value2 = $newText
    }
}
